$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58 - this shifts existing rows 58..176 down to 59..177
# (and the sheet dimension grows from A1:R176 to A1:R177 automatically).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new record.
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = "Vega Modelo de Temuco"
$ws.Range("C58").Value = "La Araucanía"
$ws.Range("D58").Value = 44791
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = 100112012
$ws.Range("G58").Value = "Espinaca"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 65
$ws.Range("K58").Value = 13000
$ws.Range("L58").Value = 13000
$ws.Range("M58").Value = 13000
$ws.Range("N58").Value = "$/docena de atados"
$ws.Range("O58").Value = "Región de La Araucanía"
$ws.Range("P58").Value = 4333
$ws.Range("Q58").Value = 3
$ws.Range("R58").Value = "Hortaliza"

# Match the date-formatted style used by the other rows' Fecha (D) column.
$ws.Range("D58").NumberFormat = $ws.Range("D59").NumberFormat
